$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2: "Oversight" text box - merge the first two runs ("Oversight " and
# "- Dr. ", using an en dash) into a single run "Oversight - Dr. " (leaves
# "Shiakolas" and " , Mr. O'Dell" runs untouched).
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shape2 = $s2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(4)
$enDash = [char]0x2013
$prefixLen = ("Oversight " + $enDash + " Dr. ").Length
$prefixRange = $para2.Characters(1, $prefixLen)
$prefixRange.Text = "Oversight " + $enDash + " Dr. "

# ---------------------------------------------------------------------------
# Slide 31: "Overall Test Results" table - fill in the Result / Pass-Fail
# columns for the five data rows.
# ---------------------------------------------------------------------------
$s31 = $p.Slides.Item(31)
$tableShape = $s31.Shapes.Item(2)
$tbl = $tableShape.Table

$rows = @(
    @{ Row = 2; Result = "100%"; PassFail = "Pass" },
    @{ Row = 3; Result = "100%"; PassFail = "Pass" },
    @{ Row = 4; Result = "100%"; PassFail = "Pass" },
    @{ Row = 5; Result = "100%"; PassFail = "Pass" },
    @{ Row = 6; Result = "30%";  PassFail = "Fail" }
)

foreach ($r in $rows) {
    $resultCell = $tbl.Cell($r.Row, 3)
    $resultCell.Shape.TextFrame.TextRange.Text = $r.Result

    $passFailCell = $tbl.Cell($r.Row, 4)
    $passFailCell.Shape.TextFrame.TextRange.Text = $r.PassFail
}
